$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-32: update Price (D) and Volume(1h) (E) columns.
# Price cells are stored as text in the source data (some contain multiple
# "." thousands separators), so we force a text number format before
# assigning the value and then restore the default style so no stray
# style index is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.431.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08095"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.218"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.389"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.358"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.651.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.757"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.457.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.908"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.26%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.216"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.832.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.915"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "

# Rows 33-34: WEMIXTOKEN and FraxShare swap positions (with updated D/E values)
$ws.Range("B33").Value = "FraxShare"
$ws.Range("C33").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.100"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "

# Rows 35-51: update Price (D) and Volume(1h) (E) columns
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02720"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08734"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2422"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06754"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6861"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6374"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.246"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.917"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.147"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.65%  "
